$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2749.25
$ws.Range("I137").Value = 1798.5
$ws.Range("J137").Value = 4333.8335
$ws.Range("K137").Value = 5395.5
$ws.Range("L137").Value = 13001.5005
$ws.Range("M137").Value = -2845.5
$ws.Range("N137").Value = -18101.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51636.6
$ws.Range("I2").Value = 1811.1
$ws.Range("J2").Value = 101462.1
$ws.Range("K2").Value = 1811.1
$ws.Range("L2").Value = 101462.1
$ws.Range("M2").Value = -1698.1
$ws.Range("N2").Value = -101688.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1830.5555
$ws.Range("I61").Value = 1171.4445
$ws.Range("J61").Value = 3148.7778
$ws.Range("K61").Value = 1171.4445
$ws.Range("L61").Value = 3148.7778
$ws.Range("M61").Value = -959.4445000000001
$ws.Range("N61").Value = -3572.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3714.2307
$ws.Range("I74").Value = 1812.8572
$ws.Range("K74").Value = 1812.8572
$ws.Range("M74").Value = -938.8571999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3714.2307
$ws.Range("I77").Value = 1812.8572
$ws.Range("K77").Value = 9064.286
$ws.Range("M77").Value = -4696.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 127697.375
$ws.Range("I102").Value = 202395.8
$ws.Range("J102").Value = 3200
$ws.Range("K102").Value = 202395.8
$ws.Range("L102").Value = 3200
$ws.Range("M102").Value = -200773.8
$ws.Range("N102").Value = -6444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 51636.6
$ws.Range("I116").Value = 1811.1
$ws.Range("J116").Value = 101462.1
$ws.Range("K116").Value = 1811.1
$ws.Range("L116").Value = 101462.1
$ws.Range("M116").Value = 482.9000000000001
$ws.Range("N116").Value = -106050.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3320.75
$ws.Range("J122").Value = 3023.8
$ws.Range("L122").Value = 9071.400000000001
$ws.Range("N122").Value = -13971.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3265.422
$ws.Range("I132").Value = 3190.923
$ws.Range("J132").Value = 3749.6667
$ws.Range("K132").Value = 9572.769
$ws.Range("L132").Value = 11249.0001
$ws.Range("M132").Value = -7042.769
$ws.Range("N132").Value = -16309.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1830.5555
$ws.Range("I136").Value = 1171.4445
$ws.Range("J136").Value = 3148.7778
$ws.Range("K136").Value = 3514.3335
$ws.Range("L136").Value = 9446.3334
$ws.Range("M136").Value = -964.3335000000002
$ws.Range("N136").Value = -14546.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51636.6
$ws.Range("I3").Value = 1811.1
$ws.Range("J3").Value = 101462.1
$ws.Range("K3").Value = 1811.1
$ws.Range("L3").Value = 101462.1
$ws.Range("M3").Value = -1697.1
$ws.Range("N3").Value = -101690.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1215.2632
$ws.Range("I99").Value = 900
$ws.Range("J99").Value = 2397.5
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 2397.5
$ws.Range("M99").Value = 598
$ws.Range("N99").Value = -5393.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 201940.1
$ws.Range("I105").Value = 127148.75
$ws.Range("K105").Value = 127148.75
$ws.Range("M105").Value = -125401.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2261.6382
$ws.Range("I134").Value = 2194.2896
$ws.Range("J134").Value = 2546
$ws.Range("K134").Value = 6582.8688
$ws.Range("L134").Value = 7638
$ws.Range("M134").Value = -4047.8688
$ws.Range("N134").Value = -12708

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 629.6896400000001
$ws.Range("I107").Value = 697.1429000000001
$ws.Range("J107").Value = 566.73334
$ws.Range("K107").Value = 697.1429000000001
$ws.Range("L107").Value = 566.73334
$ws.Range("M107").Value = 1222.8571
$ws.Range("N107").Value = -4406.73334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3526.2
$ws.Range("I122").Value = 3521.375
$ws.Range("K122").Value = 10564.125
$ws.Range("M122").Value = -8114.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4149.4116
$ws.Range("I132").Value = 3945.6667
$ws.Range("K132").Value = 11837.0001
$ws.Range("M132").Value = -9307.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 860.1818
$ws.Range("I134").Value = 844
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 2532
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = 3
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 510
$ws.Range("I60").Value = 412
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 1236
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -985
$ws.Range("N60").Value = -3502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 824.23
$ws.Range("J131").Value = 836.0833
$ws.Range("L131").Value = 2508.2499
$ws.Range("N131").Value = -12588.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2254.6667
$ws.Range("I80").Value = 2670
$ws.Range("J80").Value = 1922.4
$ws.Range("K80").Value = 2670
$ws.Range("L80").Value = 1922.4
$ws.Range("M80").Value = -1672
$ws.Range("N80").Value = -3918.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2254.6667
$ws.Range("I83").Value = 2670
$ws.Range("J83").Value = 1922.4
$ws.Range("K83").Value = 13350
$ws.Range("L83").Value = 9612
$ws.Range("M83").Value = -8358
$ws.Range("N83").Value = -19596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 531939.6
$ws.Range("I107").Value = 300.8
$ws.Range("J107").Value = 1122649.5
$ws.Range("K107").Value = 300.8
$ws.Range("L107").Value = 1122649.5
$ws.Range("M107").Value = 1619.2
$ws.Range("N107").Value = -1126489.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2951.926
$ws.Range("I122").Value = 3513.6
$ws.Range("J122").Value = 2249.8333
$ws.Range("K122").Value = 10540.8
$ws.Range("L122").Value = 6749.499899999999
$ws.Range("M122").Value = -8090.799999999999
$ws.Range("N122").Value = -11649.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2549.3044
$ws.Range("I132").Value = 2247.8462
$ws.Range("J132").Value = 2941.2
$ws.Range("K132").Value = 6743.5386
$ws.Range("L132").Value = 8823.599999999999
$ws.Range("M132").Value = -4213.5386
$ws.Range("N132").Value = -13883.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3586.0334
$ws.Range("I132").Value = 4089.9546
$ws.Range("J132").Value = 2200.25
$ws.Range("K132").Value = 12269.8638
$ws.Range("L132").Value = 6600.75
$ws.Range("M132").Value = -9739.863799999999
$ws.Range("N132").Value = -11660.75
